$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 87776.78999999999
$ws.Range("J38").Value = 184888.89
$ws.Range("L38").Value = 554666.67
$ws.Range("N38").Value = -555410.67
$ws.Range("H40").Value = 22737950
$ws.Range("J40").Value = 45465216
$ws.Range("L40").Value = 45465216
$ws.Range("N40").Value = -45465566
$ws.Range("H74").Value = 8969
$ws.Range("I74").Value = 5448
$ws.Range("K74").Value = 5448
$ws.Range("M74").Value = -4512
$ws.Range("H77").Value = 8969
$ws.Range("I77").Value = 5448
$ws.Range("K77").Value = 27240
$ws.Range("M77").Value = -22560
$ws.Range("H87").Value = 79899.8
$ws.Range("J87").Value = 79899.8
$ws.Range("L87").Value = 79899.8
$ws.Range("N87").Value = -82395.8
$ws.Range("H90").Value = 79899.8
$ws.Range("J90").Value = 79899.8
$ws.Range("L90").Value = 239699.4
$ws.Range("N90").Value = -252179.4
$ws.Range("H138").Value = 7497.08
$ws.Range("I138").Value = 6579.8335
$ws.Range("J138").Value = 7555.6274
$ws.Range("K138").Value = 19739.5005
$ws.Range("L138").Value = 22666.8822
$ws.Range("M138").Value = -14599.5005
$ws.Range("N138").Value = -32946.8822

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 17858614
$ws.Range("I74").Value = 20834478
$ws.Range("K74").Value = 20834478
$ws.Range("M74").Value = -20833604
$ws.Range("H77").Value = 17858614
$ws.Range("I77").Value = 20834478
$ws.Range("K77").Value = 104172390
$ws.Range("M77").Value = -104168022

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 1199
$ws.Range("I25").Value = 1199
$ws.Range("K25").Value = 1199
$ws.Range("M25").Value = -964
$ws.Range("H76").Value = 15761.429
$ws.Range("J76").Value = 15761.429
$ws.Range("L76").Value = 15761.429
$ws.Range("N76").Value = -16391.429
$ws.Range("H79").Value = 15761.429
$ws.Range("J79").Value = 15761.429
$ws.Range("L79").Value = 15761.429
$ws.Range("N79").Value = -17945.429
$ws.Range("H82").Value = 20346.5
$ws.Range("J82").Value = 33855.332
$ws.Range("L82").Value = 33855.332
$ws.Range("N82").Value = -34621.332
$ws.Range("H85").Value = 20346.5
$ws.Range("J85").Value = 33855.332
$ws.Range("L85").Value = 33855.332
$ws.Range("N85").Value = -36507.332
$ws.Range("H97").Value = 50000
$ws.Range("H132").Value = 118000
$ws.Range("J132").Value = 118000
$ws.Range("L132").Value = 118000
$ws.Range("N132").Value = -128120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4472.609
$ws.Range("J31").Value = 10066.333
$ws.Range("L31").Value = 10066.333
$ws.Range("N31").Value = -10656.333
$ws.Range("H34").Value = 4472.609
$ws.Range("J34").Value = 10066.333
$ws.Range("L34").Value = 10066.333
$ws.Range("N34").Value = -10470.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3530
$ws.Range("J122").Value = 6642.75
$ws.Range("L122").Value = 59784.75
$ws.Range("N122").Value = -64684.75
$ws.Range("H136").Value = 6834.5854
$ws.Range("I136").Value = 2311.625
$ws.Range("J136").Value = 9729.280000000001
$ws.Range("K136").Value = 6934.875
$ws.Range("L136").Value = 29187.84
$ws.Range("M136").Value = -1834.875
$ws.Range("N136").Value = -39387.84
$ws.Range("H137").Value = 53575500
$ws.Range("J137").Value = 10999.5
$ws.Range("L137").Value = 32998.5
$ws.Range("N137").Value = -43198.5
$ws.Range("H138").Value = 4435.5835
$ws.Range("I138").Value = 4435.5835
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 13306.7505
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -8166.750499999998
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 2617.0908
$ws.Range("I139").Value = 2617.0908
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 7851.2724
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -2711.2724
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 3996
$ws.Range("I140").Value = 3996
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 11988
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -6808
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 5095.75
$ws.Range("J141").Value = 10000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 16125.625
$ws.Range("J24").Value = 10333
$ws.Range("L24").Value = 10333
$ws.Range("N24").Value = -10679
$ws.Range("H128").Value = 95000
$ws.Range("J128").Value = 95000
$ws.Range("L128").Value = 95000
$ws.Range("N128").Value = -104960
$ws.Range("H132").Value = 7716.5454
$ws.Range("I132").Value = 7592.7646
$ws.Range("J132").Value = 8137.4
$ws.Range("K132").Value = 22778.2938
$ws.Range("L132").Value = 24412.2
$ws.Range("M132").Value = -20248.2938
$ws.Range("N132").Value = -29472.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 962.5
$ws.Range("I22").Value = 850
$ws.Range("K22").Value = 850
$ws.Range("M22").Value = -555
$ws.Range("H27").Value = 962.5
$ws.Range("I27").Value = 850
$ws.Range("K27").Value = 850
$ws.Range("M27").Value = -743
$ws.Range("H40").Value = 29170264
$ws.Range("I40").Value = 16669664
$ws.Range("K40").Value = 16669664
$ws.Range("M40").Value = -16669528
$ws.Range("H46").Value = 8916.833000000001
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 8700.200000000001
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 8700.200000000001
$ws.Range("M46").Value = -9812
$ws.Range("N46").Value = -9076.200000000001
$ws.Range("H50").Value = 27749.5
$ws.Range("J50").Value = 27749.5
$ws.Range("L50").Value = 27749.5
$ws.Range("N50").Value = -29023.5
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H68").Value = 2591.6428
$ws.Range("I68").Value = 2207.4546
$ws.Range("J68").Value = 4000.3333
$ws.Range("K68").Value = 2207.4546
$ws.Range("L68").Value = 4000.3333
$ws.Range("M68").Value = -1458.4546
$ws.Range("N68").Value = -5498.3333
$ws.Range("H71").Value = 2591.6428
$ws.Range("I71").Value = 2207.4546
$ws.Range("J71").Value = 4000.3333
$ws.Range("K71").Value = 11037.273
$ws.Range("L71").Value = 20001.6665
$ws.Range("M71").Value = -7293.273000000001
$ws.Range("N71").Value = -27489.6665
$ws.Range("H132").Value = 6822.533
$ws.Range("J132").Value = 7181.1
$ws.Range("L132").Value = 21543.3
$ws.Range("N132").Value = -26603.3
$ws.Range("H133").Value = 90955.37
$ws.Range("J133").Value = 90955.37
$ws.Range("L133").Value = 90955.37
$ws.Range("N133").Value = -96015.37
$ws.Range("H136").Value = 6297.857
$ws.Range("I136").Value = 5230
$ws.Range("K136").Value = 15690
$ws.Range("M136").Value = -13140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2838.6
$ws.Range("I96").Value = 2838.6
$ws.Range("K96").Value = 2838.6
$ws.Range("M96").Value = -1465.6
$ws.Range("H133").Value = 48669.75
$ws.Range("J133").Value = 48669.75
$ws.Range("L133").Value = 48669.75
$ws.Range("N133").Value = -58789.75
$ws.Range("H136").Value = 10407.417
$ws.Range("I136").Value = 7237.6665
$ws.Range("K136").Value = 21712.9995
$ws.Range("M136").Value = -19162.9995
